$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(311, 2, 392, 446, 13, 0, 2180, 23478, 10, 40, 2909)
    3  = @(267, 0, 31, 32, 0, 1, 0, 6129, 0, 0, 0)
    4  = @(188, 0, 33, 37, 3, 1, 833, 2424, 0, 0, 0)
    5  = @(1550, 0, 27, 29, 1, 0, 3333, 24333, 0, 0, 0)
    6  = @(1842, 0, 84, 86, 1, 0, 20, 87777, 5, 0, 0)
    7  = @(135, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    8  = @(813, 2, 234, 270, 12, 0, 727, 14734, 4, 23, 1292)
    9  = @(266, 0, 26, 27, 1, 0, 370, 4231, 0, 0, 0)
    10 = @(254, 24, 135, 197, 35, 3, 3825, 7816, 3, 0, 0)
    11 = @(1044, 1, 62, 66, 2, 1, 4, 6330, 0, 0, 0)
    12 = @(1401, 0, 60, 78, 18, 0, 8800, 14576, 0, 0, 0)
    13 = @(868, 0, 125, 126, 0, 0, 0, 11892, 2, 1, 1429)
    14 = @(268, 0, 49, 50, 1, 0, 2, 5714, 0, 0, 0)
    15 = @(484, 0, 14, 17, 3, 0, 25, 11111, 0, 0, 0)
    16 = @(191, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    17 = @(1810, 0, 19, 21, 2, 0, 20, 12857, 0, 0, 0)
    18 = @(4046, 4, 112, 268, 17, 0, 1775, 22264, 0, 135, 14456)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
